$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 61.16466974292273
$ws.Range("C2").Value = 95.32095076252179
$ws.Range("D2").Value = 99.33198101278093
$ws.Range("E2").Value = 98.9410739353898
$ws.Range("F2").Value = 98.36641915483548
$ws.Range("G2").Value = 97.84536995593169
$ws.Range("H2").Value = 97.43275980945741
$ws.Range("I2").Value = 96.1103415418012
$ws.Range("B3").Value = 70.09954864178495
$ws.Range("C3").Value = 95.33643122342397
$ws.Range("D3").Value = 99.80891427234715
$ws.Range("E3").Value = 99.04041014452028
$ws.Range("F3").Value = 98.57670418429653
$ws.Range("G3").Value = 98.08518095971731
$ws.Range("H3").Value = 97.51259823626285
$ws.Range("I3").Value = 96.09440289849485
$ws.Range("B4").Value = 82.55319384392784
$ws.Range("C4").Value = 94.94129172069931
$ws.Range("D4").Value = 99.20977570786191
$ws.Range("E4").Value = 98.71359991868967
$ws.Range("F4").Value = 98.42050176372155
$ws.Range("G4").Value = 97.93718671144732
$ws.Range("H4").Value = 97.42213464670056
$ws.Range("I4").Value = 96.05625824293379
$ws.Range("B5").Value = 75.92609584225865
$ws.Range("C5").Value = 94.90525769058259
$ws.Range("D5").Value = 99.2741154759275
$ws.Range("E5").Value = 98.88837490278451
$ws.Range("F5").Value = 98.35800185041073
$ws.Range("G5").Value = 97.86701923410335
$ws.Range("H5").Value = 97.37982321101978
$ws.Range("I5").Value = 96.07588669457674
$ws.Range("B6").Value = 72.34624428390028
$ws.Range("C6").Value = 95.39417255371639
$ws.Range("D6").Value = 99.34392425145158
$ws.Range("E6").Value = 98.86957989896591
$ws.Range("F6").Value = 98.33594835176346
$ws.Range("G6").Value = 97.97507394015781
$ws.Range("H6").Value = 97.40352571802525
$ws.Range("I6").Value = 96.02152195018505
